# Apply the "Added capability for ELN file" change to the Main_Data sheet.
#
# Summary of the change:
#   - Column A header "Plate_Well_Id" becomes "Unique_Id"
#   - Column J header "Unnamed: 13" becomes "Plate_Well_Id"
#   - For every data row, the original (un-prefixed) well id that used to
#     live in column A (e.g. "P1-A01") is copied into column J, and column A
#     is rewritten with an "SSF00607-" prefix (e.g. "SSF00607-P1-A01").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main_Data")

# Determine the last used row on the sheet.
$lastRow = $ws.UsedRange.Rows.Count

# --- Header row (row 1) ---
$ws.Cells.Item(1, 1).Value = "Unique_Id"
$ws.Cells.Item(1, 10).Value = "Plate_Well_Id"

# --- Data rows ---
for ($r = 2; $r -le $lastRow; $r++) {
    $oldId = $ws.Cells.Item($r, 1).Value2
    if ($oldId -ne $null -and $oldId -ne "") {
        $ws.Cells.Item($r, 10).Value = $oldId
        $ws.Cells.Item($r, 1).Value = "SSF00607-" + $oldId
    }
}
